$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LG Monitor (Amazon)
$ws.Range("B2").Value = 1049
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1049
$ws.Range("F2").Value = "2025-03-30 18:57:10"

# Row 3 - AOC Monitor (Amazon)
$ws.Range("B3").Value = 1004
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 1004.27
$ws.Range("F3").Value = "2025-03-30 18:57:10"

# Row 4 - iPhone (Mercado Livre)
$ws.Range("D4").Value = 4299
$ws.Range("F4").Value = "2025-03-30 18:57:11"

# Row 5 - Samsung Galaxy (Mercado Livre)
$ws.Range("D5").Value = 5158
$ws.Range("F5").Value = "2025-03-30 18:57:12"
